$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 10:39:14"
$ws1.Cells.Item(3,1).Value = "Total filas: 135"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 10:39:14"
$ws2.Cells.Item(3,1).Value = "Total filas: 15"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 10:39:14"
$ws3.Cells.Item(3,1).Value = "Total filas: 21"

$rows1 = @(
  @(33, "06:52:41", "07:23", "10_OLMOS", 31, "LP1912"),
  @(34, "07:23:38", "07:23", "16_SANTA ANA", 0, "LP1912"),
  @(35, "07:23:38", "07:31", "16_SANTA ANA", 8, "LP1912"),
  @(36, "07:23:38", "07:31", "11_ETCHEVERRY", 8, "LP1912"),
  @(77, "08:54:41", "09:22", "16_SANTA ANA", 28, "LP1912"),
  @(78, "07:59:51", "09:22", "17_ROMERO", 83, "LP1912"),
  @(79, "08:54:41", "09:23", "11_ETCHEVERRY", 29, "LP1912"),
  @(80, "08:31:16", "09:23", "16_SANTA ANA", 52, "LP1912"),
  @(104, "10:39:14", "10:40", "14_ABASTO", 1, "LP1912"),
  @(105, "09:32:47", "10:41", "17_ROMERO", 69, "LP1912"),
  @(106, "10:39:14", "10:42", "17_ROMERO", 3, "LP1912"),
  @(107, "09:32:47", "10:43", "14_ABASTO", 71, "LP1912"),
  @(108, "10:39:14", "10:47", "16_SANTA ANA", 8, "LP1912"),
  @(109, "10:39:14", "10:52", "15_ABASTO", 13, "LP1912"),
  @(110, "10:39:14", "10:53", "10_OLMOS", 14, "LP1912"),
  @(111, "09:32:47", "10:53", "27_EL RETIRO", 81, "LP1912"),
  @(112, "10:39:14", "10:57", "16_SANTA ANA", 18, "LP1912"),
  @(113, "10:39:14", "10:57", "27_EL RETIRO", 18, "LP1912"),
  @(114, "09:32:47", "11:01", "215C_EL PATO", 89, "LP1912"),
  @(115, "10:39:14", "11:02", "215C_EL PATO", 23, "LP1912"),
  @(116, "10:39:14", "11:04", "11_ETCHEVERRY", 25, "LP1912"),
  @(117, "10:39:14", "11:05", "23_HERNANDEZ", 26, "LP1912"),
  @(118, "09:32:47", "11:06", "16_P MOR-167 Y 521", 94, "LP1912"),
  @(119, "10:39:14", "11:07", "16_P MOR-167 Y 521", 28, "LP1912"),
  @(120, "10:39:14", "11:11", "10_OLMOS", 32, "LP1912"),
  @(121, "10:39:14", "11:12", "15_ABASTO", 33, "LP1912"),
  @(122, "09:32:47", "11:19", "86_EST CHICA-ESC AGRARIA", 107, "LP1912"),
  @(123, "10:39:14", "11:20", "86_EST CHICA-ESC AGRARIA", 41, "LP1912"),
  @(124, "10:39:14", "11:21", "26_HERNANDEZ", 42, "LP1912"),
  @(125, "09:32:47", "11:26", "16_P MOR-SANTA ANA", 114, "LP1912"),
  @(126, "10:39:14", "11:27", "225_C ROCA-H SUR", 48, "LP1912"),
  @(127, "10:39:14", "11:32", "81_EL PELIGRO", 53, "LP1912"),
  @(128, "10:39:14", "11:35", "23_HERNANDEZ", 56, "LP1912"),
  @(129, "10:39:14", "11:36", "11_ETCHEVERRY", 57, "LP1912"),
  @(130, "10:39:14", "11:42", "17_ROMERO", 63, "LP1912"),
  @(131, "10:39:14", "11:43", "10_OLMOS", 64, "LP1912"),
  @(132, "10:39:14", "11:51", "215B_EL PATO", 72, "LP1912"),
  @(133, "10:39:14", "11:59", "225_GOMEZ", 80, "LP1912"),
  @(134, "10:39:14", "12:02", "84_COLONIA URQUIZA-ESC 49", 83, "LP1912"),
  @(135, "10:39:14", "12:07", "16_P MOR-SANTA ANA", 88, "LP1912"),
  @(136, "10:39:14", "12:14", "17_ROMERO", 95, "LP1912"),
  @(137, "10:39:14", "12:21", "215A_EL PATO", 102, "LP1912"),
  @(138, "10:39:14", "12:21", "26_HERNANDEZ", 102, "LP1912"),
  @(139, "10:39:14", "12:22", "14_ABASTO", 103, "LP1912"),
  @(140, "10:39:14", "12:37", "27_EL RETIRO", 118, "LP1912")
)

$rows2 = @(
  @(18, "10:39:14", "11:02", "215C_EL PATO", 23, "LP1912"),
  @(19, "10:39:14", "11:51", "215B_EL PATO", 72, "LP1912"),
  @(20, "10:39:14", "12:21", "215A_EL PATO", 102, "LP1912")
)

$rows3 = @(
  @(23, "10:39:14", "10:54", "215A_LA PLATA", 15, "L6173"),
  @(25, "10:39:14", "11:14", "215C_LA PLATA", 35, "L6203"),
  @(26, "10:39:14", "12:04", "215A_LA PLATA", 85, "L6173")
)

foreach ($row in $rows1) {
  $r = $row[0]
  $ws1.Cells.Item($r,1).Value = $row[1]
  $ws1.Cells.Item($r,2).Value = $row[2]
  $ws1.Cells.Item($r,3).Value = $row[3]
  $ws1.Cells.Item($r,4).Value = $row[4]
  $ws1.Cells.Item($r,5).Value = $row[5]
}

foreach ($row in $rows2) {
  $r = $row[0]
  $ws2.Cells.Item($r,1).Value = $row[1]
  $ws2.Cells.Item($r,2).Value = $row[2]
  $ws2.Cells.Item($r,3).Value = $row[3]
  $ws2.Cells.Item($r,4).Value = $row[4]
  $ws2.Cells.Item($r,5).Value = $row[5]
}

foreach ($row in $rows3) {
  $r = $row[0]
  $ws3.Cells.Item($r,1).Value = $row[1]
  $ws3.Cells.Item($r,2).Value = $row[2]
  $ws3.Cells.Item($r,3).Value = $row[3]
  $ws3.Cells.Item($r,4).Value = $row[4]
  $ws3.Cells.Item($r,5).Value = $row[5]
}

Write-Host "Edit applied successfully."
